$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update workbook/sheet title to reflect the new "through" date
$wb.Sheets(1).Name = "Through 2021-09-21"

# Update the row label for September to reflect the new "through" date
$ws.Range("A10").Value = "September (through 09-21)"

# Update September row (row 10) values
$ws.Range("B10").Value = 23
$ws.Range("D10").Value = 49
$ws.Range("E10").Value = 40
$ws.Range("F10").Value = 50
$ws.Range("G10").Value = 82
$ws.Range("H10").Value = 127

# Update Total row (row 11) values
$ws.Range("B11").Value = 217
$ws.Range("D11").Value = 600
$ws.Range("E11").Value = 530
$ws.Range("F11").Value = 399
$ws.Range("G11").Value = 866
$ws.Range("H11").Value = 1197
